$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Update the Neo4j query text in cell B2 (shared string): replace the
# "Activity Code" / "Award Amount" lines with the new formulas.
$old = $ws.Range("B2").Text
$new = $old.Replace(
    "coalesce(p.lead_doc, '')AS ``Activity Code``,`ncoalesce(p.award_amount, '') AS ``Award Amount``,",
    "SUBSTRING(p.project_id, 1, 3) AS ``Activity code``,`n`"`$`" + apoc.number.format(toInteger(p.award_amount)) AS ``Award Amount``,"
)
$ws.Range("B2").Value = $new

# Update the active cell selection on the "startup" sheet from B19 to C20.
$ws.Activate()
$ws.Range("C20").Select()
